$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025175111852936
$ws.Range("D2").Value = 1.029422921386836
$ws.Range("E2").Value = 1.048405314612532
$ws.Range("F2").Value = 1.052741932897462
$ws.Range("I2").Value = 1.029864519276841
$ws.Range("J2").Value = 1.030345603158973
$ws.Range("K2").Value = 1.032236876639802
$ws.Range("L2").Value = 1.051165251228512
$ws.Range("M2").Value = 1.055489825073982
$ws.Range("N2").Value = 1.014121319705415
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026438496476971
$ws.Range("D3").Value = 1.030374154861914
$ws.Range("E3").Value = 1.049830620564165
$ws.Range("F3").Value = 1.054284392696868
$ws.Range("I3").Value = 1.030135829220114
$ws.Range("J3").Value = 1.031246589314984
$ws.Range("K3").Value = 1.03299586709611
$ws.Range("L3").Value = 1.052401009594728
$ws.Range("M3").Value = 1.056843320839045
$ws.Range("N3").Value = 1.014421201693478
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027254428334239
$ws.Range("D4").Value = 1.030987874395556
$ws.Range("E4").Value = 1.050752683621198
$ws.Range("F4").Value = 1.055282138236555
$ws.Range("I4").Value = 1.030308816274788
$ws.Range("J4").Value = 1.031827583154632
$ws.Range("K4").Value = 1.03348458067783
$ws.Range("L4").Value = 1.053199903994133
$ws.Range("M4").Value = 1.057718298842143
$ws.Range("N4").Value = 1.014614507050394
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027597076255585
$ws.Range("D5").Value = 1.031245455250083
$ws.Range("E5").Value = 1.051140275114702
$ws.Range("F5").Value = 1.055701517473032
$ws.Range("I5").Value = 1.030380925900057
$ws.Range("J5").Value = 1.032071355872408
$ws.Range("K5").Value = 1.033689462432309
$ws.Range("L5").Value = 1.053535590690269
$ws.Range("M5").Value = 1.058085947398079
$ws.Range("N5").Value = 1.014695596745568
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027654586784403
$ws.Range("D6").Value = 1.031288679242357
$ws.Range("E6").Value = 1.051205351016485
$ws.Range("F6").Value = 1.055771928980083
$ws.Range("I6").Value = 1.030392997433266
$ws.Range("J6").Value = 1.032112258487541
$ws.Range("K6").Value = 1.033723829413874
$ws.Range("L6").Value = 1.053591944221076
$ws.Range("M6").Value = 1.058147666086562
$ws.Range("N6").Value = 1.014709201774187
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027259008260526
$ws.Range("D7").Value = 1.03099131787806
$ws.Range("E7").Value = 1.050757862803281
$ws.Range("F7").Value = 1.055287742282052
$ws.Range("I7").Value = 1.030309782219309
$ws.Range("J7").Value = 1.031830842328882
$ws.Range("K7").Value = 1.033487320568556
$ws.Range("L7").Value = 1.053204390112383
$ws.Range("M7").Value = 1.057723212128993
$ws.Range("N7").Value = 1.014615591264868
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025602404678034
$ws.Range("D8").Value = 1.029744767602509
$ws.Range("E8").Value = 1.048887049122742
$ws.Range("F8").Value = 1.053263286045988
$ws.Range("I8").Value = 1.029956742080542
$ws.Range("J8").Value = 1.030650512364598
$ws.Range("K8").Value = 1.03249388004702
$ws.Range("L8").Value = 1.051583034332475
$ws.Range("M8").Value = 1.055947419071144
$ws.Range("N8").Value = 1.014222819640745
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022671099842422
$ws.Range("D9").Value = 1.02753436055666
$ws.Range("E9").Value = 1.045588626406615
$ws.Range("F9").Value = 1.04969315464452
$ws.Range("I9").Value = 1.029314930566821
$ws.Range("J9").Value = 1.028555142472654
$ws.Range("K9").Value = 1.030724803921773
$ws.Range("L9").Value = 1.048720238879235
$ws.Range("M9").Value = 1.052811701216056
$ws.Range("N9").Value = 1.0135250109445
$ws.Range("B10").Value = 1.019999999999999
$ws.Range("C10").Value = 1.020708451757472
$ws.Range("D10").Value = 1.026051314718886
$ws.Range("E10").Value = 1.043388125986603
$ws.Range("F10").Value = 1.047310831176735
$ws.Range("I10").Value = 1.028873749822494
$ws.Range("J10").Value = 1.027147643033881
$ws.Range("K10").Value = 1.029532839464444
$ws.Range("L10").Value = 1.046807533687228
$ws.Range("M10").Value = 1.050716498346534
$ws.Range("N10").Value = 1.013055918907177
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019856539112751
$ws.Range("D11").Value = 1.025406866519804
$ws.Range("E11").Value = 1.042434841862193
$ws.Range("F11").Value = 1.046278646403079
$ws.Range("I11").Value = 1.02867954454766
$ws.Range("J11").Value = 1.026535628208507
$ws.Range("K11").Value = 1.029013690003131
$ws.Range("L11").Value = 1.0459782531718
$ws.Range("M11").Value = 1.049808058931733
$ws.Range("N11").Value = 1.012851861984631
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01953978395638
$ws.Range("D12").Value = 1.025167144236803
$ws.Range("E12").Value = 1.042080675115745
$ws.Range("F12").Value = 1.045895146273349
$ws.Range("I12").Value = 1.028606930408409
$ws.Range("J12").Value = 1.026307910627149
$ws.Range("K12").Value = 1.028820398027003
$ws.Range("L12").Value = 1.045670055316588
$ws.Range("M12").Value = 1.049470437022487
$ws.Range("N12").Value = 1.012775924211996
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019607743439091
$ws.Range("D13").Value = 1.025218581165535
$ws.Range("E13").Value = 1.042156648565204
$ws.Range("F13").Value = 1.045977413027298
$ws.Range("I13").Value = 1.0286225280226
$ws.Range("J13").Value = 1.026356774439716
$ws.Range("K13").Value = 1.028861880527676
$ws.Range("L13").Value = 1.045736172416596
$ws.Range("M13").Value = 1.049542866635699
$ws.Range("N13").Value = 1.012792219564377
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.019830362528229
$ws.Range("D14").Value = 1.025387058067565
$ws.Range("E14").Value = 1.042405567887732
$ws.Range("F14").Value = 1.046246948259448
$ws.Range("I14").Value = 1.028673551991445
$ws.Range("J14").Value = 1.026516812939184
$ws.Range("K14").Value = 1.028997721762488
$ws.Range("L14").Value = 1.045952780864641
$ws.Range("M14").Value = 1.049780154852369
$ws.Range("N14").Value = 1.012845587847161
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.019967483345438
$ws.Range("D15").Value = 1.025490816369601
$ws.Range("E15").Value = 1.04255892521938
$ws.Range("F15").Value = 1.046413004247158
$ws.Range("I15").Value = 1.028704926211734
$ws.Range("J15").Value = 1.026615366412959
$ws.Range("K15").Value = 1.029081357419716
$ws.Range("L15").Value = 1.046086218284756
$ws.Range("M15").Value = 1.049926330975491
$ws.Range("N15").Value = 1.012878450962936
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.020764946255856
$ws.Range("D16").Value = 1.026094036336357
$ws.Range("E16").Value = 1.043451382215278
$ws.Range("F16").Value = 1.047379320145874
$ws.Range("I16").Value = 1.028886571670991
$ws.Range("J16").Value = 1.027188206229662
$ws.Range("K16").Value = 1.029567229810641
$ws.Range("L16").Value = 1.046862547329229
$ws.Range("M16").Value = 1.050776762595983
$ws.Range("N16").Value = 1.013069441652008
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021264614844283
$ws.Range("D17").Value = 1.026471807841127
$ws.Range("E17").Value = 1.044011071051686
$ws.Range("F17").Value = 1.047985292753961
$ws.Range("I17").Value = 1.028999663212896
$ws.Range("J17").Value = 1.027546846007983
$ws.Range("K17").Value = 1.029871193898925
$ws.Range("L17").Value = 1.047349228321615
$ws.Range("M17").Value = 1.051309889394816
$ws.Range("N17").Value = 1.013188993384224
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021555863426054
$ws.Range("D18").Value = 1.026691935832961
$ws.Range("E18").Value = 1.044337484824372
$ws.Range("F18").Value = 1.048338686700527
$ws.Range("I18").Value = 1.029065321715722
$ws.Range("J18").Value = 1.027755787937231
$ws.Range("K18").Value = 1.030048199725512
$ws.Range("L18").Value = 1.047632998606176
$ws.Range("M18").Value = 1.051620737622035
$ws.Range("N18").Value = 1.013258635562947
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021655137964012
$ws.Range("D19").Value = 1.026766956609732
$ws.Range("E19").Value = 1.044448776345199
$ws.Range("F19").Value = 1.048459175005839
$ws.Range("I19").Value = 1.029087657722467
$ws.Range("J19").Value = 1.027826990020439
$ws.Range("K19").Value = 1.030108504817594
$ws.Range("L19").Value = 1.04772973973142
$ws.Range("M19").Value = 1.051726709415733
$ws.Range("N19").Value = 1.013282366460905
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021211025825204
$ws.Range("D20").Value = 1.026431299285316
$ws.Range("E20").Value = 1.043951026259464
$ws.Range("F20").Value = 1.04792028384861
$ws.Range("I20").Value = 1.028987561208791
$ws.Range("J20").Value = 1.027508392887148
$ws.Range("K20").Value = 1.029838611579623
$ws.Range("L20").Value = 1.04729702270378
$ws.Range("M20").Value = 1.051252701940987
$ws.Range("N20").Value = 1.013176175966317
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019764815545961
$ws.Range("D21").Value = 1.025337455382053
$ws.Range("E21").Value = 1.042332269462764
$ws.Range("F21").Value = 1.046167579690326
$ws.Range("I21").Value = 1.028658539899929
$ws.Range("J21").Value = 1.02646969633104
$ws.Range("K21").Value = 1.028957732536502
$ws.Range("L21").Value = 1.045888999709677
$ws.Range("M21").Value = 1.049710284601968
$ws.Range("N21").Value = 1.012829876145963
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018853689722261
$ws.Range("D22").Value = 1.024647711647068
$ws.Range("E22").Value = 1.041314058069556
$ws.Range("F22").Value = 1.045064998557245
$ws.Range("I22").Value = 1.028448906963433
$ws.Range("J22").Value = 1.025814379737669
$ws.Range("K22").Value = 1.028401244842841
$ws.Range("L22").Value = 1.045002755533505
$ws.Range("M22").Value = 1.048739420135453
$ws.Range("N22").Value = 1.012611321856964
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019336870613549
$ws.Range("D23").Value = 1.025013548519271
$ws.Range("E23").Value = 1.04185387459364
$ws.Range("F23").Value = 1.045649555665852
$ws.Range("I23").Value = 1.028560299778746
$ws.Range("J23").Value = 1.026161989673807
$ws.Range("K23").Value = 1.028696501098549
$ws.Range("L23").Value = 1.045472663720596
$ws.Range("M23").Value = 1.049254198783223
$ws.Range("N23").Value = 1.0127272599182
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021235241012756
$ws.Range("D24").Value = 1.02644960403683
$ws.Range("E24").Value = 1.043978158049673
$ws.Range("F24").Value = 1.047949658758383
$ws.Range("I24").Value = 1.028993030528586
$ws.Range("J24").Value = 1.027525768958947
$ws.Range("K24").Value = 1.029853335027092
$ws.Range("L24").Value = 1.047320612491618
$ws.Range("M24").Value = 1.051278542846548
$ws.Range("N24").Value = 1.013181967884771
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023430380861075
$ws.Range("D25").Value = 1.028107457032986
$ws.Range("E25").Value = 1.046441596608538
$ws.Range("F25").Value = 1.050616486294153
$ws.Range("I25").Value = 1.029483195473406
$ws.Range("J25").Value = 1.029098698378194
$ws.Range("K25").Value = 1.031184359674261
$ws.Range("L25").Value = 1.049461053690061
$ws.Range("M25").Value = 1.053623168200728
$ws.Range("N25").Value = 1.013706091951021
